$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values updated
$ws.Range("B3").Value = 0.9991255098771341
$ws.Range("C3").Value = 0.9989474040679845
$ws.Range("D3").Value = 0.9776414451205874

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9986173871658394
$ws.Range("C4").Value = 0.9983369384392237
$ws.Range("D4").Value = 0.9868659028293542

# Row 5: AdaBoostRegressor -> MLPRegressor, values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9996270326120303
$ws.Range("C5").Value = 0.9996200389525587
$ws.Range("D5").Value = 0.9992421193657409
